# Criando uma funçao para executar o programa dentro da Classe e prints
# para mostrar a linha e qual planilha os dados foram salvos

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Salvar-Leitura {
    param($Linha, $DataHora, $Temperatura, $Umidade)

    $ws.Cells.Item($Linha, 1).Value = $DataHora
    $ws.Cells.Item($Linha, 2).Value = $Temperatura
    $ws.Cells.Item($Linha, 3).Value = $Umidade

    Write-Host "Linha $Linha salva na planilha '$($ws.Name)'"
}

Salvar-Leitura 6 "2024-11-21 22:06:49" 20 "Alerta Amarelo, Chuvas Intensas"
Salvar-Leitura 7 "2024-11-21 22:08:04" 19 "Alerta Amarelo, Chuvas Intensas"
